$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 206, shifting existing rows 206:258 down to 207:259
$ws.Rows("206:206").Insert()

# Populate the newly inserted row 206 with the new price-report record
$ws.Range("A206").Value = 5
$ws.Range("B206").Value = "Macroferia Regional de Talca"
$ws.Range("C206").Value = "Maule"
$ws.Range("D206").Value = 45093
$ws.Range("E206").Value = 7
$ws.Range("F206").Value = 100112031
$ws.Range("G206").Value = "Poroto verde"
$ws.Range("H206").Value = "Sin especificar"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 300
$ws.Range("K206").Value = 23000
$ws.Range("L206").Value = 23000
$ws.Range("M206").Value = 23000
$ws.Range("N206").Value = "$/malla 25 kilos"
$ws.Range("O206").Value = "Región de Arica y Parinacota"
$ws.Range("P206").Value = 920
$ws.Range("Q206").Value = 25
$ws.Range("R206").Value = "Hortaliza"

# Apply the date style used by the rest of the Fecha column (D) to the new cell
$ws.Range("D206").NumberFormat = $ws.Range("D207").NumberFormat
